$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a rolling weekly price log: the newest report is always
# inserted at row 2, pushing every existing data row down by one.
# Shift rows 2..68 down to rows 3..69 (copy bottom-up so data is not
# clobbered before it is copied), preserving values and formatting.
for ($r = 68; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":T" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":T" + ($r + 1))
    $src.Copy($dst)
}

# Row 2 still holds the old top entry's data (now duplicated into row 3
# as well); overwrite it in place with this week's new report values.
$ws.Range("D2").Value = 44812
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 36000
$ws.Range("O2").Value = 36000
$ws.Range("P2").Value = 36000
$ws.Range("S2").Value = 2000
